# Generate Report for Archive
# This script reorders the per-file rows in the Overview / zh-cn / de-de
# worksheets so the "5af91c6c-..." file moves from the last "in progress"
# slot up to directly after "089cdb6e-...", and marks that file's
# handback as complete (adds Latest Target File / Latest Handback File /
# Latest Handback DateTime values + hyperlinks) on the language sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Hyperlink target URLs (unchanged from the original workbook - only
# their row/column position changes).
# ---------------------------------------------------------------------
$mdUrls = @{
    "089cdb6e-55a3-440a-8cd1-d21a170433a7.md" = "https://github.com/OpenLocalizationTest/oltest/blob/4ab77f901cf1dd2aa2091bc62a82b60118c558a7/e2e/089cdb6e-55a3-440a-8cd1-d21a170433a7.md";
    "86f42771-06db-4da6-93e3-4a9e101966cb.md" = "https://github.com/OpenLocalizationTest/oltest/blob/5c40a4b705c0f23a6986d454cf1884faa56eda1c/e2e/86f42771-06db-4da6-93e3-4a9e101966cb.md";
    "fb659db4-70f5-4538-8936-7ccd74a12800.md" = "https://github.com/OpenLocalizationTest/oltest/blob/4ab77f901cf1dd2aa2091bc62a82b60118c558a7/e2e/fb659db4-70f5-4538-8936-7ccd74a12800.md";
    "5af91c6c-42ef-492e-aaa6-34afeefe9d79.md" = "https://github.com/OpenLocalizationTest/oltest/blob/7b10a5140e24d1d7194eff43e897ae3348ee294c/e2e/5af91c6c-42ef-492e-aaa6-34afeefe9d79.md";
    "e895b03e-a365-48c4-96fa-0f40707d4535.md" = "https://github.com/OpenLocalizationTest/oltest/blob/86f392ee0c1d31a2ccd89656ae64fdb57ad9ca45/e2e/e895b03e-a365-48c4-96fa-0f40707d4535.md";
}

$zhUrls = @{
    "089cdb6e-55a3-440a-8cd1-d21a170433a7.201c112032c95366c52babb0dc9cfcc40985a08e.zh-cn.xlf" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cbe88e1ad501ffaa1846f7d74661da9e023b477d/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/089cdb6e-55a3-440a-8cd1-d21a170433a7.201c112032c95366c52babb0dc9cfcc40985a08e.zh-cn.xlf";
    "86f42771-06db-4da6-93e3-4a9e101966cb.3ca688be53f46e1fac1d7c40e4b6d84b3c7b8f04.zh-cn.xlf" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c5e7070c10a8a04e54d096dc22708634f1031261/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/86f42771-06db-4da6-93e3-4a9e101966cb.3ca688be53f46e1fac1d7c40e4b6d84b3c7b8f04.zh-cn.xlf";
    "fb659db4-70f5-4538-8936-7ccd74a12800.7a592e1fd0fc852671a53e7f3b03df8e95b02793.zh-cn.xlf" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cbe88e1ad501ffaa1846f7d74661da9e023b477d/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/fb659db4-70f5-4538-8936-7ccd74a12800.7a592e1fd0fc852671a53e7f3b03df8e95b02793.zh-cn.xlf";
    "5af91c6c-42ef-492e-aaa6-34afeefe9d79.85f9ce95e4cc5e88a4e1a4411c4c3cdca349e274.zh-cn.xlf" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a4715c41302d2adc2769ffb2764a62ebb167dfdd/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/5af91c6c-42ef-492e-aaa6-34afeefe9d79.85f9ce95e4cc5e88a4e1a4411c4c3cdca349e274.zh-cn.xlf";
    "e895b03e-a365-48c4-96fa-0f40707d4535.e6936a2989f83f270f005457b64e86d158ae1d9c.zh-cn.xlf" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d8d5cbb0ab5c7f2fa849a05c7f3235508c947791/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/e895b03e-a365-48c4-96fa-0f40707d4535.e6936a2989f83f270f005457b64e86d158ae1d9c.zh-cn.xlf";
}

$deUrls = @{
    "089cdb6e-55a3-440a-8cd1-d21a170433a7.201c112032c95366c52babb0dc9cfcc40985a08e.de-de.xlf" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4f89a20149dbb21c933847b6b744dec8fc01b0ca/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/089cdb6e-55a3-440a-8cd1-d21a170433a7.201c112032c95366c52babb0dc9cfcc40985a08e.de-de.xlf";
    "86f42771-06db-4da6-93e3-4a9e101966cb.3ca688be53f46e1fac1d7c40e4b6d84b3c7b8f04.de-de.xlf" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/491923a98e122542686a9366a52ed89d9de0c39b/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/86f42771-06db-4da6-93e3-4a9e101966cb.3ca688be53f46e1fac1d7c40e4b6d84b3c7b8f04.de-de.xlf";
    "fb659db4-70f5-4538-8936-7ccd74a12800.7a592e1fd0fc852671a53e7f3b03df8e95b02793.de-de.xlf" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4f89a20149dbb21c933847b6b744dec8fc01b0ca/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/fb659db4-70f5-4538-8936-7ccd74a12800.7a592e1fd0fc852671a53e7f3b03df8e95b02793.de-de.xlf";
    "5af91c6c-42ef-492e-aaa6-34afeefe9d79.85f9ce95e4cc5e88a4e1a4411c4c3cdca349e274.de-de.xlf" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/56faf46c86190acb43913ce91c863063f7f2a739/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/5af91c6c-42ef-492e-aaa6-34afeefe9d79.85f9ce95e4cc5e88a4e1a4411c4c3cdca349e274.de-de.xlf";
    "e895b03e-a365-48c4-96fa-0f40707d4535.e6936a2989f83f270f005457b64e86d158ae1d9c.de-de.xlf" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3ec5de3ddadc9548f9da53f3a7e30d243023f97c/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/e895b03e-a365-48c4-96fa-0f40707d4535.e6936a2989f83f270f005457b64e86d158ae1d9c.de-de.xlf";
}

function Clear-SheetHyperlinks($ws) {
    $existing = @()
    foreach ($h in $ws.Hyperlinks) {
        $existing += $h
    }
    foreach ($h in $existing) {
        $h.Delete()
    }
}

# ---------------------------------------------------------------------
# 1) Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

Clear-SheetHyperlinks $wsOverview

$overviewRows = @(
    @{ Row = 2; A = "089cdb6e-55a3-440a-8cd1-d21a170433a7.md"; B = "In Translation"; C = "In Translation"; D = "2016-03-22 05:02:02" },
    @{ Row = 3; A = "5af91c6c-42ef-492e-aaa6-34afeefe9d79.md"; B = "In Translation"; C = "In Translation"; D = "2016-03-22 05:06:41" },
    @{ Row = 4; A = "86f42771-06db-4da6-93e3-4a9e101966cb.md"; B = "In Translation"; C = "In Translation"; D = "2016-03-22 05:03:26" },
    @{ Row = 5; A = "fb659db4-70f5-4538-8936-7ccd74a12800.md"; B = "In Translation"; C = "In Translation"; D = "2016-03-22 05:02:02" },
    @{ Row = 6; A = "e895b03e-a365-48c4-96fa-0f40707d4535.md"; B = "Ready for handoff"; C = "Ready for handoff"; D = "2016-03-22 05:04:12" }
)

foreach ($r in $overviewRows) {
    $row = $r.Row
    $wsOverview.Range("A$row").Value = $r.A
    $wsOverview.Range("B$row").Value = $r.B
    $wsOverview.Range("C$row").Value = $r.C
    $wsOverview.Range("D$row").Value = $r.D
    $wsOverview.Hyperlinks.Add($wsOverview.Range("A$row"), $mdUrls[$r.A], [Type]::Missing, $r.A, $r.A) | Out-Null
}

# ---------------------------------------------------------------------
# 2) zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

Clear-SheetHyperlinks $wsZh

$zhRows = @(
    @{ Row = 2; A = "089cdb6e-55a3-440a-8cd1-d21a170433a7.md"; B = ".md"; C = "In Translation"; D = "089cdb6e-55a3-440a-8cd1-d21a170433a7.201c112032c95366c52babb0dc9cfcc40985a08e.zh-cn.xlf"; E = "2016-03-22 05:01:47"; HasFG = $false; H = "0001-01-01 00:00:00"; J = "Include" },
    @{ Row = 3; A = "5af91c6c-42ef-492e-aaa6-34afeefe9d79.md"; B = ".md"; C = "In Translation"; D = "5af91c6c-42ef-492e-aaa6-34afeefe9d79.85f9ce95e4cc5e88a4e1a4411c4c3cdca349e274.zh-cn.xlf"; E = "2016-03-22 05:06:33"; HasFG = $true; F = "5af91c6c-42ef-492e-aaa6-34afeefe9d79.md"; G = "5af91c6c-42ef-492e-aaa6-34afeefe9d79.85f9ce95e4cc5e88a4e1a4411c4c3cdca349e274.zh-cn.xlf"; H = "2016-03-22 05:07:13"; J = "Include" },
    @{ Row = 4; A = "86f42771-06db-4da6-93e3-4a9e101966cb.md"; B = ".md"; C = "In Translation"; D = "86f42771-06db-4da6-93e3-4a9e101966cb.3ca688be53f46e1fac1d7c40e4b6d84b3c7b8f04.zh-cn.xlf"; E = "2016-03-22 05:03:18"; HasFG = $false; H = "0001-01-01 00:00:00"; J = "Include" },
    @{ Row = 5; A = "fb659db4-70f5-4538-8936-7ccd74a12800.md"; B = ".md"; C = "In Translation"; D = "fb659db4-70f5-4538-8936-7ccd74a12800.7a592e1fd0fc852671a53e7f3b03df8e95b02793.zh-cn.xlf"; E = "2016-03-22 05:01:47"; HasFG = $false; H = "0001-01-01 00:00:00"; J = "Include" },
    @{ Row = 6; A = "e895b03e-a365-48c4-96fa-0f40707d4535.md"; B = ".md"; C = "Ready for handoff"; D = "e895b03e-a365-48c4-96fa-0f40707d4535.e6936a2989f83f270f005457b64e86d158ae1d9c.zh-cn.xlf"; E = "2016-03-22 05:04:05"; HasFG = $false; H = "0001-01-01 00:00:00"; J = "Include" }
)

foreach ($r in $zhRows) {
    $row = $r.Row
    $wsZh.Range("A$row").Value = $r.A
    $wsZh.Range("B$row").Value = $r.B
    $wsZh.Range("C$row").Value = $r.C
    $wsZh.Range("D$row").Value = $r.D
    $wsZh.Range("E$row").Value = $r.E
    $wsZh.Range("H$row").Value = $r.H
    $wsZh.Range("J$row").Value = $r.J

    $wsZh.Hyperlinks.Add($wsZh.Range("A$row"), $mdUrls[$r.A], [Type]::Missing, $r.A, $r.A) | Out-Null
    $wsZh.Hyperlinks.Add($wsZh.Range("D$row"), $zhUrls[$r.D], [Type]::Missing, $r.D, $r.D) | Out-Null

    if ($r.HasFG) {
        $wsZh.Range("F$row").Value = $r.F
        $wsZh.Range("G$row").Value = $r.G
        $wsZh.Hyperlinks.Add($wsZh.Range("F$row"), $mdUrls[$r.F], [Type]::Missing, $r.F, $r.F) | Out-Null
        $wsZh.Hyperlinks.Add($wsZh.Range("G$row"), $zhUrls[$r.G], [Type]::Missing, $r.G, $r.G) | Out-Null
    }
}

# ---------------------------------------------------------------------
# 3) de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

Clear-SheetHyperlinks $wsDe

$deRows = @(
    @{ Row = 2; A = "089cdb6e-55a3-440a-8cd1-d21a170433a7.md"; B = ".md"; C = "In Translation"; D = "089cdb6e-55a3-440a-8cd1-d21a170433a7.201c112032c95366c52babb0dc9cfcc40985a08e.de-de.xlf"; E = "2016-03-22 05:02:02"; HasFG = $false; H = "0001-01-01 00:00:00"; J = "Include" },
    @{ Row = 3; A = "5af91c6c-42ef-492e-aaa6-34afeefe9d79.md"; B = ".md"; C = "In Translation"; D = "5af91c6c-42ef-492e-aaa6-34afeefe9d79.85f9ce95e4cc5e88a4e1a4411c4c3cdca349e274.de-de.xlf"; E = "2016-03-22 05:06:41"; HasFG = $true; F = "5af91c6c-42ef-492e-aaa6-34afeefe9d79.md"; G = "5af91c6c-42ef-492e-aaa6-34afeefe9d79.85f9ce95e4cc5e88a4e1a4411c4c3cdca349e274.de-de.xlf"; H = "2016-03-22 05:07:27"; J = "Include" },
    @{ Row = 4; A = "86f42771-06db-4da6-93e3-4a9e101966cb.md"; B = ".md"; C = "In Translation"; D = "86f42771-06db-4da6-93e3-4a9e101966cb.3ca688be53f46e1fac1d7c40e4b6d84b3c7b8f04.de-de.xlf"; E = "2016-03-22 05:03:26"; HasFG = $false; H = "0001-01-01 00:00:00"; J = "Include" },
    @{ Row = 5; A = "fb659db4-70f5-4538-8936-7ccd74a12800.md"; B = ".md"; C = "In Translation"; D = "fb659db4-70f5-4538-8936-7ccd74a12800.7a592e1fd0fc852671a53e7f3b03df8e95b02793.de-de.xlf"; E = "2016-03-22 05:02:02"; HasFG = $false; H = "0001-01-01 00:00:00"; J = "Include" },
    @{ Row = 6; A = "e895b03e-a365-48c4-96fa-0f40707d4535.md"; B = ".md"; C = "Ready for handoff"; D = "e895b03e-a365-48c4-96fa-0f40707d4535.e6936a2989f83f270f005457b64e86d158ae1d9c.de-de.xlf"; E = "2016-03-22 05:04:12"; HasFG = $false; H = "0001-01-01 00:00:00"; J = "Include" }
)

foreach ($r in $deRows) {
    $row = $r.Row
    $wsDe.Range("A$row").Value = $r.A
    $wsDe.Range("B$row").Value = $r.B
    $wsDe.Range("C$row").Value = $r.C
    $wsDe.Range("D$row").Value = $r.D
    $wsDe.Range("E$row").Value = $r.E
    $wsDe.Range("H$row").Value = $r.H
    $wsDe.Range("J$row").Value = $r.J

    $wsDe.Hyperlinks.Add($wsDe.Range("A$row"), $mdUrls[$r.A], [Type]::Missing, $r.A, $r.A) | Out-Null
    $wsDe.Hyperlinks.Add($wsDe.Range("D$row"), $deUrls[$r.D], [Type]::Missing, $r.D, $r.D) | Out-Null

    if ($r.HasFG) {
        $wsDe.Range("F$row").Value = $r.F
        $wsDe.Range("G$row").Value = $r.G
        $wsDe.Hyperlinks.Add($wsDe.Range("F$row"), $mdUrls[$r.F], [Type]::Missing, $r.F, $r.F) | Out-Null
        $wsDe.Hyperlinks.Add($wsDe.Range("G$row"), $deUrls[$r.G], [Type]::Missing, $r.G, $r.G) | Out-Null
    }
}

Write-Host "Done applying Generate Report for Archive changes."
